# UC2 ekstra tilbehør + pris
# Adds a new glossary row ("ekstra tilbehør") to the Definitions sheet and
# extends the formatted-but-empty rows below it (mirrors the table being
# resized from A3:E16 to A3:E43).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Definitions")

# --- Row 17: new glossary term -------------------------------------------
# Write the new shared strings in the same order the source workbook uses
# them (Definition/info, Aliases, Term) so they line up with the rest of
# the edit.
$ws.Cells.Item(17, 2).Value = " antal af bagage, barnevogne eller autostole, samt behov for hjælpemidler - SF-UC2-bestilFlextur"
$ws.Cells.Item(17, 5).Value = "ikke-obligatorisk oplsyninger"
$ws.Cells.Item(17, 1).Value = "ekstra tilbehør"

$ws.Rows.Item(17).RowHeight = 57.6

# A/C/D/E use the centered+wrapped glossary-row style; B keeps the
# left-aligned, vertically centered + wrapped "definition" style.
foreach ($col in 1, 3, 4, 5) {
    $cell = $ws.Cells.Item(17, $col)
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4108
    $cell.WrapText = $true
}
$bCell = $ws.Cells.Item(17, 2)
$bCell.VerticalAlignment = -4108
$bCell.WrapText = $true

# --- Rows 18-43: extend the blank, formatted glossary rows ---------------
for ($r = 18; $r -le 43; $r++) {
    foreach ($col in 1, 3, 4, 5) {
        $cell = $ws.Cells.Item($r, $col)
        $cell.HorizontalAlignment = -4108
        $cell.VerticalAlignment = -4108
        $cell.WrapText = $true
    }
}

# --- Extend the Term/Definition table + autofilter to the new range ------
$tbl = $ws.ListObjects.Item("Table2")
$tbl.Resize($ws.Range("A3:E43"))

# --- Sheet view bookkeeping (matches the author's final cursor state) ----
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Application.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("B19").Select()
